$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 10248
$ws1.Range("F8").Value = 145
$ws1.Range("F11").Value = 4069
$ws1.Range("F18").Value = 2257
$ws1.Range("F22").Value = 8687
$ws1.Range("F24").Value = 1012
$ws1.Range("F30").Value = 1905
$ws1.Range("F35").Value = 24
$ws1.Range("F44").Value = 1455
$ws1.Range("F45").Value = 2341
$ws1.Range("F46").Value = 791
$ws1.Range("F48").Value = 1239

# Sheet "全部类型" (all types) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 10248
$ws4.Range("F8").Value = 145
$ws4.Range("F14").Value = 4069
$ws4.Range("F22").Value = 8687
$ws4.Range("F25").Value = 1012
$ws4.Range("F30").Value = 1905
$ws4.Range("F33").Value = 24
$ws4.Range("F40").Value = 1455
$ws4.Range("F42").Value = 2341
$ws4.Range("F43").Value = 795
$ws4.Range("F48").Value = 1239
